# "fix get all data excel change to get last once"
#
# The order-receipt log in Sheet1 (A1:A6) is updated so the last few rows
# reflect the corrected ("get last one" instead of "get all") export logic:
#   - A3: MacBook Air order is now Chip M2 / Color Silver
#   - A4: MacBook Air order is now Chip M1 / Color Midnight
#   - A5: Asus Tuf Gaming F15 order is now cpu gen 9 / ram 8GB
#   - A6: same corrected Asus Tuf Gaming F15 order as A5 (duplicate "last" write)
# A1/A2 are untouched historical rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "คุณสั่งชื้อ MacBook Air `nโดยมี `nChip: Macbook Air M2 `nColor: Silver `nในราคา: 55000"
$ws.Range("A4").Value = "คุณสั่งชื้อ MacBook Air `nโดยมี `nChip: Macbook Air M1 `nColor: Midnight `nในราคา: 55000"
$ws.Range("A5").Value = "คุณสั่งชื้อ Asus Tuf Gaming F15 `nโดยมี `ncpu: Intel core I5 gen 9 `nram: 8GB `ndisk: 512GB SSD `ngpu:NVIDIA GeForce GTX 1650 `nในราคา: 35000"
$ws.Range("A6").Value = "คุณสั่งชื้อ Asus Tuf Gaming F15 `nโดยมี `ncpu: Intel core I5 gen 9 `nram: 8GB `ndisk: 512GB SSD `ngpu:NVIDIA GeForce GTX 1650 `nในราคา: 35000"

# The live session's cursor ended up past the data (A11) when this was saved.
$ws.Range("A11").Select()
